$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$changes = @(
    @{Row=2; Col=4; Value='30.350.19'},
    @{Row=2; Col=5; Value='  -1.27%  '},
    @{Row=3; Col=4; Value='1.889.93'},
    @{Row=3; Col=5; Value='  -1.39%  '},
    @{Row=4; Col=4; Value='1.001'},
    @{Row=4; Col=5; Value='  -0.13%  '},
    @{Row=5; Col=4; Value='237.77'},
    @{Row=5; Col=5; Value='  -1.51%  '},
    @{Row=6; Col=4; Value='1.001'},
    @{Row=6; Col=5; Value='  -0.06%  '},
    @{Row=7; Col=4; Value='0.4825'},
    @{Row=7; Col=5; Value='  -1.96%  '},
    @{Row=8; Col=4; Value='0.2895'},
    @{Row=8; Col=5; Value='  -3.72%  '},
    @{Row=9; Col=4; Value='0.06596'},
    @{Row=9; Col=5; Value='  -2.81%  '},
    @{Row=10; Col=4; Value='1.889.90'},
    @{Row=10; Col=5; Value='  -1.32%  '},
    @{Row=11; Col=4; Value='16.86'},
    @{Row=11; Col=5; Value='  -2.34%  '},
    @{Row=12; Col=4; Value='0.07448'},
    @{Row=12; Col=5; Value='  +1.45%  '},
    @{Row=13; Col=4; Value='5.162'},
    @{Row=13; Col=5; Value='  -1.19%  '},
    @{Row=14; Col=4; Value='87.52'},
    @{Row=14; Col=5; Value='  -1.65%  '},
    @{Row=15; Col=4; Value='0.6617'},
    @{Row=15; Col=5; Value='  -2.26%  '},
    @{Row=16; Col=4; Value='30.338.19'},
    @{Row=17; Col=4; Value='13.39'},
    @{Row=17; Col=5; Value='  -2.05%  '},
    @{Row=18; Col=2; Value='Dai'},
    @{Row=18; Col=3; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Row=18; Col=4; Value='1.000'},
    @{Row=18; Col=5; Value='  -0.06%  '},
    @{Row=19; Col=2; Value='ShibaInu'},
    @{Row=19; Col=3; Value='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'},
    @{Row=19; Col=4; Value='0.000007756'},
    @{Row=19; Col=5; Value='  -3.13%  '},
    @{Row=20; Col=2; Value='Uniswap'},
    @{Row=20; Col=3; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'},
    @{Row=20; Col=4; Value='5.436'},
    @{Row=20; Col=5; Value='  +0.40%  '},
    @{Row=21; Col=2; Value='WrappedliquidstakedEther2.0'},
    @{Row=21; Col=3; Value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{Row=21; Col=4; Value='2.142.16'},
    @{Row=21; Col=5; Value='  -0.98%  '},
    @{Row=22; Col=4; Value='1.001'},
    @{Row=22; Col=5; Value='  -0.15%  '},
    @{Row=23; Col=4; Value='193.10'},
    @{Row=23; Col=5; Value='  -4.31%  '},
    @{Row=24; Col=4; Value='6.171'},
    @{Row=24; Col=5; Value='  -2.58%  '},
    @{Row=25; Col=4; Value='9.411'},
    @{Row=25; Col=5; Value='  -3.10%  '},
    @{Row=26; Col=4; Value='162.96'},
    @{Row=26; Col=5; Value='  +1.12%  '},
    @{Row=27; Col=4; Value='18.21'},
    @{Row=27; Col=5; Value='  -3.59%  '},
    @{Row=28; Col=4; Value='1.947'},
    @{Row=28; Col=5; Value='  -1.30%  '},
    @{Row=29; Col=5; Value='  -0.87%  '},
    @{Row=30; Col=4; Value='4.293'},
    @{Row=30; Col=5; Value='  -2.01%  '},
    @{Row=31; Col=4; Value='0.09121'},
    @{Row=31; Col=5; Value='  -0.64%  '},
    @{Row=32; Col=4; Value='4.044'},
    @{Row=32; Col=5; Value='  -1.40%  '},
    @{Row=33; Col=4; Value='0.05100'},
    @{Row=33; Col=5; Value='  -4.12%  '},
    @{Row=34; Col=4; Value='1.155'},
    @{Row=34; Col=5; Value='  +2.46%  '},
    @{Row=35; Col=4; Value='0.7325'},
    @{Row=35; Col=5; Value='  -2.24%  '},
    @{Row=36; Col=4; Value='2.711'},
    @{Row=36; Col=5; Value='  +0.35%  '},
    @{Row=37; Col=4; Value='0.01799'},
    @{Row=37; Col=5; Value='  -3.64%  '},
    @{Row=38; Col=4; Value='2.646'},
    @{Row=38; Col=5; Value='  -2.94%  '},
    @{Row=39; Col=4; Value='0.9175'},
    @{Row=39; Col=5; Value='  -1.66%  '},
    @{Row=40; Col=4; Value='2.081'},
    @{Row=40; Col=5; Value='  -0.95%  '},
    @{Row=41; Col=4; Value='5.906'},
    @{Row=41; Col=5; Value='  -1.15%  '},
    @{Row=42; Col=4; Value='106.49'},
    @{Row=42; Col=5; Value='  -1.07%  '},
    @{Row=43; Col=5; Value='  -4.27%  '},
    @{Row=44; Col=5; Value='  +0.05%  '},
    @{Row=45; Col=4; Value='7.558'},
    @{Row=45; Col=5; Value='  -2.49%  '},
    @{Row=46; Col=4; Value='0.1330'},
    @{Row=46; Col=5; Value='  -5.14%  '},
    @{Row=47; Col=4; Value='1.562'},
    @{Row=47; Col=5; Value='  +7.20%  '},
    @{Row=48; Col=4; Value='64.79'},
    @{Row=48; Col=5; Value='  -11.38%  '},
    @{Row=49; Col=4; Value='8.954'},
    @{Row=49; Col=5; Value='  -2.92%  '},
    @{Row=50; Col=5; Value='  -3.21%  '},
    @{Row=51; Col=4; Value='33.98'},
    @{Row=51; Col=5; Value='  -5.64%  '}
)

foreach ($item in $changes) {
    Set-TextValue $ws.Cells.Item($item.Row, $item.Col) $item.Value
}